$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the used range dimension implicitly by writing to AF47 as well.

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell (C1 uses style index 1) onto the new header cells.
$ws.Range("C1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-47: Wins=78, Losses=84, Ties=0 for every team row.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
